$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.826.94'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +3.88%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.495.40'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.48%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '491.90'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.78'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +10.78%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  +1.89%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.510.33'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.16%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0992'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.56%  '

$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.72'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.93%  '

$ws.Range("E12").Value = '  +2.79%  '

$ws.Range("E13").Value = '  +1.38%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.934.38'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.34%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.849.60'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +3.59%  '

$ws.Range("E16").Value = '  +3.45%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.504.29'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.68%  '

$ws.Range("E19").Value = '  +4.49%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.30'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.80%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.26'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.39%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.90'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +4.00%  '

$ws.Range("E24").Value = '  +1.89%  '

$ws.Range("E25").Value = '  +0.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.69%  '

$ws.Range("E27").Value = '  -2.20%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.609.56'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.49%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.63'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +4.42%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0811'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +4.47%  '

$ws.Range("E31").Value = '  +0.12%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '151.79'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.47%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.27'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.80%  '

$ws.Range("E34").Value = '  +3.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.29'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.97%  '

$ws.Range("E36").Value = '  +5.25%  '

$ws.Range("E37").Value = '  +4.61%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.877'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.66%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.40'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +8.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '34.23'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.52'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.82%  '

$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.617'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.66%  '

$ws.Range("B43").Value = 'Hedera'
$ws.Range("C43").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0562'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.90%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.995'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.21%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '267.90'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +5.88%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.81'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0939'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.94%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.22'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0228'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.96'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +4.51%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.892.54'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.10%  '
